# deleted duplicate fastq files in 3357 (and 3275), and corrected fastq file in 4144 01.13.20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicate first data row (row 2) - this shifts every
# subsequent row up by one, so the old row 3 becomes the new row 2, etc.,
# and the trailing row (old 29) disappears, shrinking the used range
# from A1:G29 down to A1:G28.
$ws.Rows(2).Delete()

# Leave the just-deleted row selected, matching Excel's behavior after
# a row delete (selects the full row that now occupies that position).
$ws.Rows(2).Select() | Out-Null
